$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"
$wsOverview.Range("B9").Value = "In Translation"
$wsOverview.Range("C9").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C8").Value = "In Translation"
$wsZhCn.Range("C9").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C8").Value = "In Translation"
$wsDeDe.Range("C9").Value = "In Translation"
